# Actualización automática 2025-11-12 11:30:07
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M4").Value = 1653.75
$wsVentasGrupo.Range("M14").Value = 3053.72
$wsVentasGrupo.Range("M26").Value = "3 de 24"

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F4").Value = 1653.75
$wsVentaMensual.Range("F14").Value = 3053.72
$wsVentaMensual.Range("F26").Value = 5468.83

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# ColumnWidth uses character units with an internal padding offset (~0.8333)
# relative to the raw OOXML <col width> value, so compensate to land on 13 / 24.
$wsCumplimiento.Columns.Item(4).ColumnWidth = 12.166666666666666
$wsCumplimiento.Columns.Item(6).ColumnWidth = 23.166666666666668

$wsCumplimiento.Range("D12").Value = 5025.39
$wsCumplimiento.Range("E12").Value = 29675.61
$wsCumplimiento.Range("F12").Value = 0.1448197458286505

$wsCumplimiento.Range("D14").Value = 5468.83
$wsCumplimiento.Range("E14").Value = 35308.91058948192
$wsCumplimiento.Range("F14").Value = 0.1341131195829573
